$d = $word.ActiveDocument

function Replace-InRange($rng, $old, $new) {
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: failed to find/replace: $old"
    }
    return $ok
}

function Replace-InPara($paraIndex, $old, $new) {
    $rng = $d.Paragraphs($paraIndex).Range
    return Replace-InRange $rng $old $new
}

# --- Title ---
Replace-InPara 1 "Decoding the Enigma of Consciousness" "History: The Tapestry of Time" | Out-Null

# --- Author line: "Dr" + "." + " Alex Rodriguez" -> "Alexander Thornton" ---
Replace-InPara 2 "Dr. Alex Rodriguez" "Alexander Thornton" | Out-Null

# --- Email line: "alex" / "." / "rodriguez@consciousnessstudies" / "." / "org"
#     -> "alex@schoolstaff" / "." / "net"  (i.e. alex@schoolstaff.net) ---
Replace-InPara 3 "alex.rodriguez@consciousnessstudies.org" "alex@schoolstaff.net" | Out-Null

# --- Body paragraph (paragraph 5): expand each "sentence slot" into the new sentences ---
Replace-InPara 5 "The enigma of consciousness persists as one of the most captivating and elusive mysteries in science" "History is a dynamic and ever-evolving narrative of human civilization, a vast landscape of interconnected events, people, and ideas that have shaped the world we inhabit today" | Out-Null

Replace-InPara 5 " What is consciousness, and how does it arise from the physical processes of the brain? The journey to unravel this enigma has led us through the labyrinth of neuroscience, philosophy, psychology, and beyond" " As we delve into the chronicles of history, we embark on a journey through diverse cultures, empires, and revolutions, each leaving an indelible imprint on the tapestry of time. From the ancient civilizations of Mesopotamia and Egypt to the technological advancements of the modern era, history provides invaluable insights into the human experience, offering lessons that resonate across generations" | Out-Null

Replace-InPara 5 "From the intricate neural networks of the brain to the subjective experiences of the self, consciousness remains a realm where myriad theories converge and diverge" "At its core, history is a collective story of ambition, resilience, and innovation" | Out-Null

Replace-InPara 5 " Is it an emergent property, a product of complex interactions among neurons, or does it transcend the physical realm entirely? The very nature of consciousness defies simple explanations, beckoning us to explore the depths of our own being" " We witness the rise and fall of empires, the birth of new ideologies, and the indomitable spirit of individuals who dared to challenge the status quo. From the battles fought on distant battlefields to the quiet struggles for justice, history unveils the intricate connections between past events and present realities. It teaches us about the consequences of unchecked power, the fragility of peace, and the enduring legacy of human ingenuity" | Out-Null

Replace-InPara 5 "Consciousness, with its kaleidoscope of colors, textures, emotions, and thoughts, poses challenges to our conventional understanding of the universe" "History also offers a profound understanding of the human condition" | Out-Null

Replace-InPara 5 " Can consciousness be reduced to mere electrochemical reactions, or does it hint at a deeper, more fundamental reality? As we probe the enigma of consciousness, we encounter fundamental questions about our place in the universe and the nature of reality itself" " Through the exploration of different time periods, we gain insights into the motivations and aspirations of individuals from all walks of life. We learn about the sacrifices made by ordinary people, the decisions that shaped the course of nations, and the profound impact of cultural exchange. History reminds us of our shared humanity, transcending boundaries of race, religion, and ethnicity" | Out-Null

# --- Summary heading paragraph stays "Summary" (no change) ---

# --- Summary body paragraph (paragraph 7) ---
Replace-InPara 7 "The enigma of consciousness continues to captivate and challenge our understanding of the universe" "In essence, history is a dynamic and multi-faceted discipline that encompasses a vast array of human experiences" | Out-Null

Replace-InPara 7 " Its exploration spans multiple disciplines, from neuroscience to philosophy and psychology" " By studying history, we gain invaluable insights into the past, present, and future" | Out-Null

Replace-InPara 7 " The search for answers to questions about the nature of consciousness, its relationship to the physical brain, and its implications for our understanding of reality remains an ongoing quest" " History teaches us about the interconnectedness of human societies, the challenges faced by our ancestors, and the lessons that can be learned from their triumphs and failures" | Out-Null

Replace-InPara 7 " The journey to decode the enigma of consciousness promises to illuminate the deepest mysteries of the human experience and expand our horizons of knowledge" " It helps us understand the complexity of the world we live in and provides a foundation for making informed decisions about the future" | Out-Null

# --- Add a new empty paragraph at the very end of the document (after the Summary body) ---
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
